$d = $word.ActiveDocument

# 1. Title date: "April 30, 2024" -> "May 07, 2024"
$d.Content.Find.Execute("April 30, 2024", $true, $false, $false, $false, $false, $true, 1, $false, "May 07, 2024", 2)

# 2. Italic query-instruction paragraph text
$oldQuery = "Extract any quote that includes a national action or plan that addresses " + [char]8220 + "{variable_name}" + [char]8221 + " which we define as " + [char]8220 + "{variable_description}" + [char]8221 + ". Only include direct quotation with the corresponding page number(s) with a brief explanation of the context of this quote within the text. It is very important not to hallucinate."
$newQuery = "Extract any quote mentions " + [char]8220 + "{variable_name}" + [char]8221 + ". Only include direct quotation with the corresponding page number(s). "
$d.Content.Find.Execute($oldQuery, $true, $false, $false, $false, $false, $true, 1, $false, $newQuery, 2)

# 3. First table, data row: "SDG 1" -> "Cement"
$tbl1 = $d.Tables.Item(1)
$tbl1.Cell(2,1).Range.Text = "Cement"

# 4. First table, data row: clear the "End poverty..." description cell
$tbl1.Cell(2,2).Range.InsertXML("<w:p xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'><w:r/></w:p>")

# 5. Second section heading: document filename
$d.Content.Find.Execute("Canada_Canada's Enhanced NDC Submission1_FINAL EN.pdf", $true, $false, $false, $false, $false, $true, 1, $false, "IRE03 CLIMATE ACTION PLAN 2023.pdf", 2)

# 6/7. Second table header row
$tbl2 = $d.Tables.Item(2)
$tbl2.Cell(1,1).Range.Text = "Quote"
$tbl2.Cell(1,2).Range.Text = "Related Variables"

# 8. Second table, data row: "SDG 1" -> "Cement"
$tbl2.Cell(2,1).Range.Text = "Cement"

# 9. Second table, data row: replace long single quote with many line-broken quotes
$tbl2.Cell(2,2).Range.Text = "Ireland will need to follow suit, and the public sector – collectively the State’s largest purchaser of construction projects – will lead by specifying lower carbon cement and concrete for future construction projects. [page 153 `vIndustrial process emissions from cement production are recognised as being hard to abate. [page 153 `vOptions that are available include reducing the clinker content of the final product; using alternative construction materials and methods to displace cement; and physically preventing the emissions from going into the atmosphere by capturing them and placing them in long-term storage. [page 153 `vEI will also support the high-risk research and development phase of work with cement and construction materials companies to develop novel cementitious materials; innovative products; and more efficient production technologies. [page 153 `vBy reducing the clinker content of cement through the use of novel binders and fillers we can reduce the carbon intensity of cement without compromising its integrity. [page 153 `vSupport is required to further increase the market share for timber construction and to displace the demand for cement in the construction sector. [page 153 `vActions will include a programme of work to require public bodies to specify low carbon cement products, where practicable, for public sector construction projects, and to identify suitable construction projects to assess the carbon impact of alternative construction materials through suitable whole life-cycle analysis approaches. [page 164 `vDecrease energy related emissions associated with cement production through fuel switching and efficiencies. [page 151 `vAll public bodies shall: - Cease using disposable cups, plates and cutlery from any public sector canteen or closed facility, excluding clinical (i.e., non-canteen healthcare) environments; - Specify low carbon construction methods and low carbon cement material as far as practicable for directly procured or supported construction projects from 2023. [page 108 `vProcess emissions are those generated during the manufacturing process, such as the release of CO2 from limestone during cement clinker production. [page 148 `vActively deliver a series of measures to reduce embodied carbon in construction materials, and emissions from cement production. [page 152 `vSpecify low carbon construction methods and low carbon cement material as far as practicable for directly procured or supported construction projects from 2023. [page 157 `vThe Cement Task Force shall prepare and submit to Government a public procurement policy by no later than Q2 2023 to facilitate public bodies to incorporate the principle of low carbon construction methods and materials and whole life-cycle analysis approaches in all publicly procured or supported projects. [page 113 `v"

# 10. Footer processing-summary line
$d.Content.Find.Execute("1 documents (42 total pages) processed in 6.10 seconds", $true, $false, $false, $false, $false, $true, 1, $false, "1 documents (284 total pages) processed in 55.10 seconds", 2)
